# Insert a new "roundCount" column before the old column C.
# This shifts existing columns C..N one place to the right (-> D..O)
# and refreshes the benchmark numbers in the data rows to match the
# newly uploaded results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").EntireColumn.Insert()

$ws.Range("C1").Value = "roundCount"

# Row 2
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 0.005930138
$ws.Range("G2").Value = 0.008047676
$ws.Range("H2").Value = 0.006680596
$ws.Range("I2").Value = 0.002848268
$ws.Range("J2").Value = 17236787.2
$ws.Range("K2").Value = 17236787.2
$ws.Range("L2").Value = 17242316.8
$ws.Range("M2").Value = 17287577.6

# Row 3
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 0.006154358
$ws.Range("G3").Value = 0.008313572
$ws.Range("H3").Value = 0.006794679
$ws.Range("I3").Value = 0.002935338
$ws.Range("J3").Value = 17391616
$ws.Range("K3").Value = 17391616
$ws.Range("L3").Value = 17391616
$ws.Range("M3").Value = 17391616

# Row 4
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 0.006088567
$ws.Range("G4").Value = 0.008310997
$ws.Range("H4").Value = 0.006939733
$ws.Range("I4").Value = 0.002926743
$ws.Range("J4").Value = 17391616
$ws.Range("K4").Value = 17391616
$ws.Range("L4").Value = 17391616
$ws.Range("M4").Value = 17391616

# Row 5
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 0.006283724
$ws.Range("G5").Value = 0.008261144
$ws.Range("H5").Value = 0.006801677
$ws.Range("I5").Value = 0.002909243
$ws.Range("J5").Value = 17391616
$ws.Range("K5").Value = 17391616
$ws.Range("L5").Value = 17391616
$ws.Range("M5").Value = 17391616

# Row 6
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 0.00616343
$ws.Range("G6").Value = 0.008422697
$ws.Range("H6").Value = 0.006927955
$ws.Range("I6").Value = 0.003017271
$ws.Range("J6").Value = 17391616
$ws.Range("K6").Value = 17391616
$ws.Range("L6").Value = 17391616
$ws.Range("M6").Value = 17391616

# Row 7
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 0.006241477
$ws.Range("G7").Value = 0.008450532
$ws.Range("H7").Value = 0.006804407
$ws.Range("I7").Value = 0.002949333
$ws.Range("J7").Value = 17391616
$ws.Range("K7").Value = 17391616
$ws.Range("L7").Value = 17391616
$ws.Range("M7").Value = 17391616

